$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric values in C1 and D1
$ws.Range("C1").Value = 4
$ws.Range("D1").Value = 65

# E1 must hold the literal text string "true" (not the boolean TRUE).
# A leading apostrophe forces Excel to store the entry as text instead
# of auto-converting it to a boolean; resetting the style back to
# "Normal" afterwards clears the quote-prefix formatting mark that the
# apostrophe entry leaves behind, so only the cell's value/type change.
$ws.Range("E1").Value = "'true"
$ws.Range("E1").Style = "Normal"
